# Revert "adding term 2.0 now utf-8"
# This reverts the workbook's Metadata (version/date/contact) back to the
# previous release, and restores the "Include from FSIII" sheet to its
# prior (smaller) 15-row layout.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Metadata -----------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B3").Value  = "1.1.0"                               # Version
$ws1.Range("B8").Value  = "2023-07-10T23:08:03+02:00"           # Date
$ws1.Range("B10").Value = "No display for ContactDetail"        # Contact

# --- Sheet 2: Include from FSIII -------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Drop rows 16-27 (the extra concept rows added in the forward commit),
# shrinking the sheet back to 15 rows (A1:B15).
$ws2.Range("A16:B27").EntireRow.Delete()

# Restore the 12 nursing-area concept codes in rows 2-13 (column B stays
# blank for these rows, as it already was).
$codes = @("I9", "I2", "I11", "I1", "I5", "I7", "I4", "I8", "I3", "I12", "I10", "I6")
for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = 2 + $i
    $ws2.Cells.Item($row, 1).Value = $codes[$i]
}

# Row 14 is a blank spacer row.
$ws2.Cells.Item(14, 1).Value = ""
$ws2.Cells.Item(14, 2).Value = ""

# Row 15 holds the System URI reference.
$ws2.Cells.Item(15, 1).Value = "System URI"
$ws2.Cells.Item(15, 2).Value = "urn:oid:1.2.208.176.2.21"
